$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "244.34"),
    @("E2", "-0.70%"),
    @("D3", "27.11"),
    @("E3", "3.70%"),
    @("D4", "5.157"),
    @("E4", "1.43%"),
    @("D5", "0.05630"),
    @("E5", "0.22%"),
    @("D6", "6.472"),
    @("E6", "-0.04%"),
    @("D7", "0.8164"),
    @("E7", "0.67%"),
    @("D8", "0.8308"),
    @("E8", "-2.01%"),
    @("D9", "0.1327"),
    @("E9", "-1.44%"),
    @("D10", "0.06900"),
    @("E10", "-0.93%"),
    @("D11", "0.02943"),
    @("E11", "6.42%"),
    @("D12", "0.09393"),
    @("E12", "-0.12%"),
    @("E13", "-0.03%"),
    @("D14", "0.04244"),
    @("E14", "-9.86%"),
    @("D15", "0.0005995"),
    @("E15", "-93.89%"),
    @("D16", "0.006132"),
    @("E16", "0.83%"),
    @("D17", "3.565"),
    @("E17", "0.21%"),
    @("D18", "3.017"),
    @("E18", "-0.04%"),
    @("E19", "9.01%"),
    @("D21", "0.03108"),
    @("E21", "-3.77%"),
    @("D23", "3.752"),
    @("E23", "0.14%"),
    @("E24", "-0.11%"),
    @("D25", "0.001226"),
    @("E25", "-1.68%"),
    @("D26", "0.004485"),
    @("E26", "-2.91%"),
    @("E27", "2.07%"),
    @("E28", "-0.48%"),
    @("D40", "0.03651"),
    @("E40", "-0.15%"),
    @("D41", "0.006037"),
    @("E41", "-1.16%"),
    @("D42", "0.1053"),
    @("E42", "-0.27%"),
    @("D43", "0.002299"),
    @("E43", "-8.01%"),
    @("D44", "0.008223"),
    @("E44", "-4.75%"),
    @("D45", "0.00005416"),
    @("E45", "2.37%"),
    @("E46", "-0.01%"),
    @("E47", "-18.08%"),
    @("D48", "0.002642"),
    @("E48", "28.97%"),
    @("E49", "-0.01%"),
    @("E50", "-0.01%")
)

foreach ($u in $updates) {
    $c = $ws.Range($u[0])
    $c.NumberFormat = "@"
    $c.Value = $u[1]
    $c.Style = "Normal"
}
